$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Merge "the private key in" + " " into a single run "the private key in "
# -----------------------------------------------------------------------
$d.Content.Find.Execute("the private key in ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the private key in ", 2) | Out-Null

# -----------------------------------------------------------------------
# 2) Before the "Use Chrome" paragraph (Production section), insert:
#      - text "Prepare users for downtime:   " added to that same
#        paragraph (keeping its lastRenderedPageBreak)
#      - a new paragraph with the "update app_version..." command and a
#        "_GoBack" bookmark
#      - a new empty paragraph
#      - a new paragraph containing "Use Chrome" (no page-break marker)
# -----------------------------------------------------------------------
$rngUseChrome = $d.Content
$rngUseChrome.Find.Execute("Use Chrome", $true, $false, $false, $false, $false,
                            $true, 1, $false) | Out-Null
$useChromePara = $rngUseChrome.Paragraphs(1)
$useChromeParaStart = $useChromePara.Range.Duplicate
$useChromeParaStart.Collapse(1)

$xmlPrepare = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Prepare users for downtime:   </w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:firstLine="720"/></w:pPr><w:bookmarkStart w:id="15" w:name="_GoBack"/><w:bookmarkEnd w:id="15"/><w:proofErr w:type="gramStart"/><w:r><w:t>update</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>app_version</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> set </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>downmessage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = some message about going down when/how long</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Use Chrome</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$useChromeParaStart.InsertXML($xmlPrepare)

# -----------------------------------------------------------------------
# 3) Remove the old "_GoBack" bookmark that used to sit after
#    "/webapp/gf/download" (a new one was just created above, and Word
#    only ever keeps one "_GoBack" bookmark).
# -----------------------------------------------------------------------
$rngFtp = $d.Content
$rngFtp.Find.Execute("FTP the file to ", $true, $false, $false, $false, $false,
                      $true, 1, $false) | Out-Null
$ftpPara = $rngFtp.Paragraphs(1)
$ftpParaStart = $ftpPara.Range.Duplicate
$ftpParaStart.Collapse(1)
$xmlFtp = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">FTP the file to </w:t></w:r><w:r><w:t>/webapp/gf/download</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ftpParaStart.InsertXML($xmlFtp)

# -----------------------------------------------------------------------
# 4) Move the lastRenderedPageBreak marker from the "Beta testing mostly
#    ..." run up onto the "Beta Testing" heading run.
# -----------------------------------------------------------------------
$rngBetaHeading = $d.Content
$rngBetaHeading.Find.Execute("Beta Testing", $true, $false, $false, $false, $false,
                              $true, 1, $false) | Out-Null
$betaHeadingPara = $rngBetaHeading.Paragraphs(1)
$betaHeadingParaStart = $betaHeadingPara.Range.Duplicate
$betaHeadingParaStart.Collapse(1)
$xmlBetaHeading = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="23" w:name="_Toc367458775"/><w:r><w:lastRenderedPageBreak/><w:t>Beta Testing</w:t></w:r><w:bookmarkEnd w:id="23"/></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$betaHeadingParaStart.InsertXML($xmlBetaHeading)

$rngBetaBody = $d.Content
$rngBetaBody.Find.Execute("Beta testing mostly pertains", $true, $false, $false, $false, $false,
                           $true, 1, $false) | Out-Null
$betaBodyPara = $rngBetaBody.Paragraphs(1)
$betaBodyParaStart = $betaBodyPara.Range.Duplicate
$betaBodyParaStart.Collapse(1)
$xmlBetaBody = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Beta testing mostly pertains to the mobile at this time (9/19/13)</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$betaBodyParaStart.InsertXML($xmlBetaBody)
